# Generate Report for Archive
#
# The "a10bb7dc" source file has moved on from "Ready for handoff" to
# "In Translation", while the "5e3adcf0" source file is still sitting at
# "Ready for handoff". Every report sheet (Overview, zh-cn, de-de) lists
# the two source files in rows 2-3, so on each sheet those two rows swap
# their file identity (and all of the per-file detail that goes with it)
# and the a10bb7dc row's status becomes "In Translation".

$wb = $excel.ActiveWorkbook

$md5 = "5e3adcf0-2dfd-4bd9-9041-3078e0acd852.md"
$mda = "a10bb7dc-2af9-4a10-ba96-0ca14e954cea.md"

# Blue (#6495ED), underlined -- matches this workbook's "HyperLink" look.
$hyperlinkColor = 15570276

function Restyle-Hyperlink($range) {
    $range.Font.Underline = 2
    $range.Font.Color = $hyperlinkColor
}

# ---------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = $mda
$ws.Range("B2").Value = "In Translation"
$ws.Range("C2").Value = "In Translation"

$ws.Range("A3").Value = $md5
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/904365b1a364532a90c517f26c8d0e24b6b8edaa/e2e/5e3adcf0-2dfd-4bd9-9041-3078e0acd852.md", [System.Type]::Missing, [System.Type]::Missing, $mda) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/c9a2f7b91eadab47f82d1e34dda6f13c4e1ad256/e2e/a10bb7dc-2af9-4a10-ba96-0ca14e954cea.md", [System.Type]::Missing, [System.Type]::Missing, $md5) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/c9a2f7b91eadab47f82d1e34dda6f13c4e1ad256/.localization-config", [System.Type]::Missing, [System.Type]::Missing, ".localization-config") | Out-Null

Restyle-Hyperlink $ws.Range("A2")
Restyle-Hyperlink $ws.Range("A3")
Restyle-Hyperlink $ws.Range("A4")

# ---------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = $mda
$ws.Range("B2").Value = "In Translation"
$ws.Range("C2").Value = "a10bb7dc-2af9-4a10-ba96-0ca14e954cea.2d3f83a0351f35ed02ad7ecee08b2de4891d8c6e.zh-cn.xlf"
$ws.Range("D2").Value = "2016-02-22 13:47:00"

$ws.Range("A3").Value = $md5
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "5e3adcf0-2dfd-4bd9-9041-3078e0acd852.be8a2ac0bdfac6a0c26fa9cad11a1af283b2562e.zh-cn.xlf"
$ws.Range("D3").Value = "2016-02-22 13:45:48"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/904365b1a364532a90c517f26c8d0e24b6b8edaa/e2e/5e3adcf0-2dfd-4bd9-9041-3078e0acd852.md", [System.Type]::Missing, [System.Type]::Missing, $mda) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e30bce61b2ccf4fd628866be3658b41fd2c9f242/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/5e3adcf0-2dfd-4bd9-9041-3078e0acd852.be8a2ac0bdfac6a0c26fa9cad11a1af283b2562e.zh-cn.xlf", [System.Type]::Missing, [System.Type]::Missing, "a10bb7dc-2af9-4a10-ba96-0ca14e954cea.2d3f83a0351f35ed02ad7ecee08b2de4891d8c6e.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/c9a2f7b91eadab47f82d1e34dda6f13c4e1ad256/e2e/a10bb7dc-2af9-4a10-ba96-0ca14e954cea.md", [System.Type]::Missing, [System.Type]::Missing, $md5) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9997573a4eb74395eee9bfa559aff92b94f8fb2f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/a10bb7dc-2af9-4a10-ba96-0ca14e954cea.2d3f83a0351f35ed02ad7ecee08b2de4891d8c6e.zh-cn.xlf", [System.Type]::Missing, [System.Type]::Missing, "5e3adcf0-2dfd-4bd9-9041-3078e0acd852.be8a2ac0bdfac6a0c26fa9cad11a1af283b2562e.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/c9a2f7b91eadab47f82d1e34dda6f13c4e1ad256/.localization-config", [System.Type]::Missing, [System.Type]::Missing, ".localization-config") | Out-Null

Restyle-Hyperlink $ws.Range("A2")
Restyle-Hyperlink $ws.Range("A3")
Restyle-Hyperlink $ws.Range("A4")
Restyle-Hyperlink $ws.Range("C2")
Restyle-Hyperlink $ws.Range("C3")

# ---------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = $mda
$ws.Range("B2").Value = "In Translation"
$ws.Range("C2").Value = "a10bb7dc-2af9-4a10-ba96-0ca14e954cea.2d3f83a0351f35ed02ad7ecee08b2de4891d8c6e.de-de.xlf"
$ws.Range("D2").Value = "2016-02-22 13:47:15"

$ws.Range("A3").Value = $md5
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "5e3adcf0-2dfd-4bd9-9041-3078e0acd852.be8a2ac0bdfac6a0c26fa9cad11a1af283b2562e.de-de.xlf"
$ws.Range("D3").Value = "2016-02-22 13:46:04"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/904365b1a364532a90c517f26c8d0e24b6b8edaa/e2e/5e3adcf0-2dfd-4bd9-9041-3078e0acd852.md", [System.Type]::Missing, [System.Type]::Missing, $mda) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3b61d5a74457a7aa614eb0ad167fdc561f96b011/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/5e3adcf0-2dfd-4bd9-9041-3078e0acd852.be8a2ac0bdfac6a0c26fa9cad11a1af283b2562e.de-de.xlf", [System.Type]::Missing, [System.Type]::Missing, "a10bb7dc-2af9-4a10-ba96-0ca14e954cea.2d3f83a0351f35ed02ad7ecee08b2de4891d8c6e.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/c9a2f7b91eadab47f82d1e34dda6f13c4e1ad256/e2e/a10bb7dc-2af9-4a10-ba96-0ca14e954cea.md", [System.Type]::Missing, [System.Type]::Missing, $md5) | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/62e544b02ebcbea8d9edfc5196c633ba78a91f1b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/a10bb7dc-2af9-4a10-ba96-0ca14e954cea.2d3f83a0351f35ed02ad7ecee08b2de4891d8c6e.de-de.xlf", [System.Type]::Missing, [System.Type]::Missing, "5e3adcf0-2dfd-4bd9-9041-3078e0acd852.be8a2ac0bdfac6a0c26fa9cad11a1af283b2562e.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/c9a2f7b91eadab47f82d1e34dda6f13c4e1ad256/.localization-config", [System.Type]::Missing, [System.Type]::Missing, ".localization-config") | Out-Null

Restyle-Hyperlink $ws.Range("A2")
Restyle-Hyperlink $ws.Range("A3")
Restyle-Hyperlink $ws.Range("A4")
Restyle-Hyperlink $ws.Range("C2")
Restyle-Hyperlink $ws.Range("C3")
